# Add a new "2024/11/30" data column (CE) to the 合成確率 sheet,
# one day after the existing last column (CD, 2024/11/29).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new column CE (83) the same width as the other data columns.
# (ColumnWidth 11.17 is what Excel reports back for the existing "12"-stored
# columns, e.g. CD, so re-using it keeps the stored width identical.)
$ws.Columns("CE").ColumnWidth = 11.17

# Row 1 header: the new date. Prefix with an apostrophe so Excel keeps it as
# literal text "2024/11/30" instead of auto-converting it to a date serial,
# then copy CD1's format (font/fill) onto it so it matches the other headers.
$ws.Range("CE1").Value = "'2024/11/30"
$ws.Range("CD1").Copy()
$ws.Range("CE1").PasteSpecial(-4122)

# Data rows 2-53: copy the format of an existing cell in the same row that
# already carries the desired highlight style, then write the new value.
$ws.Range("N2").Copy()
$ws.Range("CE2").PasteSpecial(-4122)
$ws.Range("CE2").Value = 129
$ws.Range("A3").Copy()
$ws.Range("CE3").PasteSpecial(-4122)
$ws.Range("CE3").Value = 163.9
$ws.Range("A4").Copy()
$ws.Range("CE4").PasteSpecial(-4122)
$ws.Range("CE4").Value = 199.6
$ws.Range("O5").Copy()
$ws.Range("CE5").PasteSpecial(-4122)
$ws.Range("CE5").Value = 135.6
$ws.Range("A6").Copy()
$ws.Range("CE6").PasteSpecial(-4122)
$ws.Range("CE6").Value = 148.4
$ws.Range("A7").Copy()
$ws.Range("CE7").PasteSpecial(-4122)
$ws.Range("CE7").Value = 173.6
$ws.Range("B8").Copy()
$ws.Range("CE8").PasteSpecial(-4122)
$ws.Range("CE8").Value = 131.6
$ws.Range("F9").Copy()
$ws.Range("CE9").PasteSpecial(-4122)
$ws.Range("CE9").Value = 117.7
$ws.Range("J10").Copy()
$ws.Range("CE10").PasteSpecial(-4122)
$ws.Range("CE10").Value = 113.7
$ws.Range("S11").Copy()
$ws.Range("CE11").PasteSpecial(-4122)
$ws.Range("CE11").Value = 132.9
$ws.Range("C12").Copy()
$ws.Range("CE12").PasteSpecial(-4122)
$ws.Range("CE12").Value = 121.4
$ws.Range("D13").Copy()
$ws.Range("CE13").PasteSpecial(-4122)
$ws.Range("CE13").Value = 130.4
$ws.Range("A14").Copy()
$ws.Range("CE14").PasteSpecial(-4122)
$ws.Range("CE14").Value = 151.8
$ws.Range("H15").Copy()
$ws.Range("CE15").PasteSpecial(-4122)
$ws.Range("CE15").Value = 132.3
$ws.Range("L16").Copy()
$ws.Range("CE16").PasteSpecial(-4122)
$ws.Range("CE16").Value = 115
$ws.Range("C17").Copy()
$ws.Range("CE17").PasteSpecial(-4122)
$ws.Range("CE17").Value = 126.8
$ws.Range("A18").Copy()
$ws.Range("CE18").PasteSpecial(-4122)
$ws.Range("CE18").Value = 161.5
$ws.Range("H19").Copy()
$ws.Range("CE19").PasteSpecial(-4122)
$ws.Range("CE19").Value = 138.8
$ws.Range("H20").Copy()
$ws.Range("CE20").PasteSpecial(-4122)
$ws.Range("CE20").Value = 137.5
$ws.Range("A21").Copy()
$ws.Range("CE21").PasteSpecial(-4122)
$ws.Range("CE21").Value = 143.6
$ws.Range("A22").Copy()
$ws.Range("CE22").PasteSpecial(-4122)
$ws.Range("CE22").Value = 231
$ws.Range("R23").Copy()
$ws.Range("CE23").PasteSpecial(-4122)
$ws.Range("CE23").Value = 122.9
$ws.Range("A24").Copy()
$ws.Range("CE24").PasteSpecial(-4122)
$ws.Range("CE24").Value = 173
$ws.Range("E25").Copy()
$ws.Range("CE25").PasteSpecial(-4122)
$ws.Range("CE25").Value = 109.8
$ws.Range("A26").Copy()
$ws.Range("CE26").PasteSpecial(-4122)
$ws.Range("CE26").Value = 169.5
$ws.Range("A27").Copy()
$ws.Range("CE27").PasteSpecial(-4122)
$ws.Range("CE27").Value = 140
$ws.Range("B28").Copy()
$ws.Range("CE28").PasteSpecial(-4122)
$ws.Range("CE28").Value = 129.1
$ws.Range("A29").Copy()
$ws.Range("CE29").PasteSpecial(-4122)
$ws.Range("CE29").Value = 161.3
$ws.Range("E30").Copy()
$ws.Range("CE30").PasteSpecial(-4122)
$ws.Range("CE30").Value = 139.4
$ws.Range("L31").Copy()
$ws.Range("CE31").PasteSpecial(-4122)
$ws.Range("CE31").Value = 126.2
$ws.Range("A32").Copy()
$ws.Range("CE32").PasteSpecial(-4122)
$ws.Range("CE32").Value = 210
$ws.Range("O33").Copy()
$ws.Range("CE33").PasteSpecial(-4122)
$ws.Range("CE33").Value = 118.3
$ws.Range("C34").Copy()
$ws.Range("CE34").PasteSpecial(-4122)
$ws.Range("CE34").Value = 125.5
$ws.Range("A35").Copy()
$ws.Range("CE35").PasteSpecial(-4122)
$ws.Range("CE35").Value = 171.4
$ws.Range("A36").Copy()
$ws.Range("CE36").PasteSpecial(-4122)
$ws.Range("CE36").Value = 151
$ws.Range("A37").Copy()
$ws.Range("CE37").PasteSpecial(-4122)
$ws.Range("CE37").Value = 152.5
$ws.Range("A38").Copy()
$ws.Range("CE38").PasteSpecial(-4122)
$ws.Range("CE38").Value = 172.4
$ws.Range("A39").Copy()
$ws.Range("CE39").PasteSpecial(-4122)
$ws.Range("CE39").Value = 148.3
$ws.Range("A40").Copy()
$ws.Range("CE40").PasteSpecial(-4122)
$ws.Range("CE40").Value = 179.2
$ws.Range("A41").Copy()
$ws.Range("CE41").PasteSpecial(-4122)
$ws.Range("CE41").Value = 157.4
$ws.Range("A42").Copy()
$ws.Range("CE42").PasteSpecial(-4122)
$ws.Range("CE42").Value = 164.6
$ws.Range("A43").Copy()
$ws.Range("CE43").PasteSpecial(-4122)
$ws.Range("CE43").Value = 143.8
$ws.Range("B44").Copy()
$ws.Range("CE44").PasteSpecial(-4122)
$ws.Range("CE44").Value = 132.1
$ws.Range("C45").Copy()
$ws.Range("CE45").PasteSpecial(-4122)
$ws.Range("CE45").Value = 133
$ws.Range("A46").Copy()
$ws.Range("CE46").PasteSpecial(-4122)
$ws.Range("CE46").Value = 170.6
$ws.Range("A47").Copy()
$ws.Range("CE47").PasteSpecial(-4122)
$ws.Range("CE47").Value = 160.1
$ws.Range("F48").Copy()
$ws.Range("CE48").PasteSpecial(-4122)
$ws.Range("CE48").Value = 120.9
$ws.Range("A49").Copy()
$ws.Range("CE49").PasteSpecial(-4122)
$ws.Range("CE49").Value = 178.9
$ws.Range("V50").Copy()
$ws.Range("CE50").PasteSpecial(-4122)
$ws.Range("CE50").Value = 126.5
$ws.Range("A51").Copy()
$ws.Range("CE51").PasteSpecial(-4122)
$ws.Range("CE51").Value = 182.3
$ws.Range("A52").Copy()
$ws.Range("CE52").PasteSpecial(-4122)
$ws.Range("CE52").Value = 141.6
$ws.Range("A53").Copy()
$ws.Range("CE53").PasteSpecial(-4122)
$ws.Range("CE53").Value = 167.9

$excel.CutCopyMode = 0
Write-Host "Added 2024/11/30 column (CE) with 53 rows"
